$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for data rows 2..387 is updated from 45205 to 45206
# (serial date 2023-10-06 -> 2023-10-07) for every row.
$ws.Range("C2:C387").Value = 45206
